# Scheduled runner: refresh Universalis market-price derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) per leve row.
# Values below mirror the latest market snapshot across all 8 crafting-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 729.65
$ws.Range("J33").Value = 1900
$ws.Range("L33").Value = 1900
$ws.Range("N33").Value = -2358

$ws.Range("H40").Value = 2999

$ws.Range("H64").Value = 60008.844
$ws.Range("I64").Value = 6227
$ws.Range("K64").Value = 6227
$ws.Range("M64").Value = -5979

$ws.Range("H67").Value = 60008.844
$ws.Range("I67").Value = 6227
$ws.Range("K67").Value = 6227
$ws.Range("M67").Value = -5369

$ws.Range("H98").Value = 2003
$ws.Range("I98").Value = 2099.1333
$ws.Range("K98").Value = 2099.1333
$ws.Range("M98").Value = -601.1333

$ws.Range("H122").Value = 2003
$ws.Range("I122").Value = 2099.1333
$ws.Range("K122").Value = 6297.3999
$ws.Range("M122").Value = -3847.3999

$ws.Range("H132").Value = 225515.48
$ws.Range("I132").Value = 240268.8
$ws.Range("J132").Value = 4215.8
$ws.Range("K132").Value = 720806.3999999999
$ws.Range("L132").Value = 12647.4
$ws.Range("M132").Value = -718276.3999999999
$ws.Range("N132").Value = -17707.4

$ws.Range("H138").Value = 3208.6155
$ws.Range("J138").Value = 3728.7576
$ws.Range("L138").Value = 11186.2728
$ws.Range("N138").Value = -21466.2728


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4908480
$ws.Range("I32").Value = 5470046
$ws.Range("J32").Value = 14834.429
$ws.Range("K32").Value = 5470046
$ws.Range("L32").Value = 14834.429
$ws.Range("M32").Value = -5469759
$ws.Range("N32").Value = -15408.429

$ws.Range("H61").Value = 7229.613
$ws.Range("I61").Value = 2803.1333
$ws.Range("K61").Value = 2803.1333
$ws.Range("M61").Value = -2591.1333

$ws.Range("H136").Value = 7229.613
$ws.Range("I136").Value = 2803.1333
$ws.Range("K136").Value = 8409.3999
$ws.Range("M136").Value = -5859.3999

$ws.Range("H137").Value = 88000
$ws.Range("J137").Value = 88000
$ws.Range("L137").Value = 88000
$ws.Range("N137").Value = -98200


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3556.9524
$ws.Range("I20").Value = 3240.1333
$ws.Range("J20").Value = 4349
$ws.Range("K20").Value = 3240.1333
$ws.Range("L20").Value = 4349
$ws.Range("M20").Value = -2993.1333
$ws.Range("N20").Value = -4843

$ws.Range("H134").Value = 18255
$ws.Range("I134").Value = 20836.736
$ws.Range("J134").Value = 7729.4614
$ws.Range("K134").Value = 62510.208
$ws.Range("L134").Value = 23188.3842
$ws.Range("M134").Value = -59975.208
$ws.Range("N134").Value = -28258.3842


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 25142
$ws.Range("J9").Value = 25142
$ws.Range("L9").Value = 25142
$ws.Range("N9").Value = -25478

$ws.Range("H31").Value = 7111.4644
$ws.Range("I31").Value = 2157.9375
$ws.Range("J31").Value = 9092.875
$ws.Range("K31").Value = 2157.9375
$ws.Range("L31").Value = 9092.875
$ws.Range("M31").Value = -1862.9375
$ws.Range("N31").Value = -9682.875

$ws.Range("H34").Value = 7111.4644
$ws.Range("I34").Value = 2157.9375
$ws.Range("J34").Value = 9092.875
$ws.Range("K34").Value = 2157.9375
$ws.Range("L34").Value = 9092.875
$ws.Range("M34").Value = -1955.9375
$ws.Range("N34").Value = -9496.875

$ws.Range("H58").Value = 3728.8235
$ws.Range("I58").Value = 3199.6155
$ws.Range("K58").Value = 3199.6155
$ws.Range("M58").Value = -2996.6155

$ws.Range("H136").Value = 3728.8235
$ws.Range("I136").Value = 3199.6155
$ws.Range("K136").Value = 9598.8465
$ws.Range("M136").Value = -7048.8465

$ws.Range("H138").Value = 54925.4
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws.Range("H141").Value = 281466.72
$ws.Range("J141").Value = 320853.16
$ws.Range("L141").Value = 320853.16
$ws.Range("N141").Value = -331213.16


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2969.9614
$ws.Range("I5").Value = 436.5
$ws.Range("J5").Value = 5141.5
$ws.Range("K5").Value = 1309.5
$ws.Range("L5").Value = 15424.5
$ws.Range("M5").Value = -1197.5
$ws.Range("N5").Value = -15648.5

$ws.Range("H7").Value = 5730048.5
$ws.Range("I7").Value = 6666779.5
$ws.Range("J7").Value = 5027500.5
$ws.Range("K7").Value = 20000338.5
$ws.Range("L7").Value = 15082501.5
$ws.Range("M7").Value = -20000226.5
$ws.Range("N7").Value = -15082725.5

$ws.Range("H34").Value = 1971.1428
$ws.Range("I34").Value = 933
$ws.Range("J34").Value = 2749.75
$ws.Range("K34").Value = 2799
$ws.Range("L34").Value = 8249.25
$ws.Range("M34").Value = -2715
$ws.Range("N34").Value = -8417.25

$ws.Range("H92").Value = 632.55
$ws.Range("I92").Value = 716.3570999999999
$ws.Range("K92").Value = 2149.0713
$ws.Range("M92").Value = -901.0712999999996

$ws.Range("H107").Value = 505.76923
$ws.Range("J107").Value = 487.5
$ws.Range("L107").Value = 1462.5
$ws.Range("N107").Value = -5302.5

$ws.Range("H135").Value = 2969.9614
$ws.Range("I135").Value = 436.5
$ws.Range("J135").Value = 5141.5
$ws.Range("K135").Value = 3928.5
$ws.Range("L135").Value = 46273.5
$ws.Range("M135").Value = -1393.5
$ws.Range("N135").Value = -51343.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5012.4546
$ws.Range("I122").Value = 2455.4375
$ws.Range("K122").Value = 7366.3125
$ws.Range("M122").Value = -4916.3125

$ws.Range("H132").Value = 2546.9678
$ws.Range("I132").Value = 2121.5
$ws.Range("J132").Value = 3320.5454
$ws.Range("K132").Value = 6364.5
$ws.Range("L132").Value = 9961.636200000001
$ws.Range("M132").Value = -3834.5
$ws.Range("N132").Value = -15021.6362


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3000
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 3000
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = -2830
$ws.Range("N12").Value = -3340

$ws.Range("H40").Value = 4099.82
$ws.Range("I40").Value = 5340.75
$ws.Range("K40").Value = 5340.75
$ws.Range("M40").Value = -5204.75

$ws.Range("H141").Value = 97956
$ws.Range("J141").Value = 97956
$ws.Range("L141").Value = 97956
$ws.Range("N141").Value = -108316


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 15000
$ws.Range("I3").Value = 15000
$ws.Range("K3").Value = 15000
$ws.Range("M3").Value = -14886

$ws.Range("H100").Value = 1693.8
$ws.Range("I100").Value = 1054.1578
$ws.Range("J100").Value = 2085.8386
$ws.Range("K100").Value = 2108.3156
$ws.Range("L100").Value = 4171.6772
$ws.Range("M100").Value = -1567.3156
$ws.Range("N100").Value = -5253.6772

$ws.Range("H122").Value = 4010.0557
$ws.Range("I122").Value = 3716.5293
$ws.Range("K122").Value = 11149.5879
$ws.Range("M122").Value = -8699.5879

$ws.Range("H136").Value = 10188548
$ws.Range("I136").Value = 13322449
$ws.Range("J136").Value = 3367.875
$ws.Range("K136").Value = 39967347
$ws.Range("L136").Value = 10103.625
$ws.Range("M136").Value = -39964797
$ws.Range("N136").Value = -15203.625

